$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "243.55"
Set-TextValue $ws "D3" "23.69"
Set-TextValue $ws "B4" "HuobiToken"
Set-TextValue $ws "C4" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws "D4" "5.298"
Set-TextValue $ws "E4" "3HuobiTokenHT"
Set-TextValue $ws "B5" "Cronos"
Set-TextValue $ws "C5" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D5" "0.05794"
Set-TextValue $ws "E5" "4CronosCRO"
Set-TextValue $ws "B6" "KuCoinToken"
Set-TextValue $ws "C6" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws "D6" "6.496"
Set-TextValue $ws "E6" "5KuCoinTokenKCS"
Set-TextValue $ws "B7" "GateToken"
Set-TextValue $ws "C7" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws "D7" "3.343"
Set-TextValue $ws "E7" "6GateTokenGT"
Set-TextValue $ws "B8" "MXToken"
Set-TextValue $ws "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D8" "0.8086"
Set-TextValue $ws "E8" "7MXTokenMX"
Set-TextValue $ws "B9" "FTXToken"
Set-TextValue $ws "C9" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws "D9" "0.8774"
Set-TextValue $ws "E9" "8FTXTokenFTT"
Set-TextValue $ws "B10" "WazirX"
Set-TextValue $ws "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws "D10" "0.1388"
Set-TextValue $ws "E10" "9WazirXWRX"
Set-TextValue $ws "B11" "MandalaExchangeToken"
Set-TextValue $ws "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws "D11" "0.07278"
Set-TextValue $ws "E11" "10MandalaExchangeTokenMDX"
Set-TextValue $ws "B12" "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws "D12" "0.03073"
Set-TextValue $ws "E12" "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue $ws "B13" "BitrueCoin"
Set-TextValue $ws "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D13" "0.03055"
Set-TextValue $ws "E13" "12BitrueCoinBTR"
Set-TextValue $ws "B14" "BitMartToken"
Set-TextValue $ws "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws "D14" "0.09324"
Set-TextValue $ws "E14" "13BitMartTokenBMX"
Set-TextValue $ws "B15" "MCDex"
Set-TextValue $ws "C15" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws "D15" "3.860"
Set-TextValue $ws "E15" "14MCDexMCB"
Set-TextValue $ws "B16" "BitForexToken"
Set-TextValue $ws "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D16" "0.001543"
Set-TextValue $ws "E16" "15BitForexTokenBF"
Set-TextValue $ws "B17" "CoinExToken"
Set-TextValue $ws "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws "D17" "0.04718"
Set-TextValue $ws "E17" "16CoinExTokenCET"
Set-TextValue $ws "B18" "One"
Set-TextValue $ws "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D18" "0.0006013"
Set-TextValue $ws "E18" "17OneONE"
Set-TextValue $ws "B19" "TigerCash"
Set-TextValue $ws "C19" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D19" "0.006096"
Set-TextValue $ws "E19" "18TigerCashTCH"
Set-TextValue $ws "B20" "BitKan"
Set-TextValue $ws "C20" "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws "D20" "0.001268"
Set-TextValue $ws "E20" "19BitKanKAN"
Set-TextValue $ws "B21" "HotbitToken"
Set-TextValue $ws "C21" "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws "D21" "0.004593"
Set-TextValue $ws "E21" "20HotbitTokenHTB"
Set-TextValue $ws "B22" "NitroEx"
Set-TextValue $ws "C22" "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws "D22" "0.00008702"
Set-TextValue $ws "E22" "21NitroExNTX"
Set-TextValue $ws "B23" "LEO"
Set-TextValue $ws "C23" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D23" "3.578"
Set-TextValue $ws "E23" "22LEOLEO"
Set-TextValue $ws "D25" "0.3211"
Set-TextValue $ws "D28" "0.0002345"
Set-TextValue $ws "D40" "0.03780"
Set-TextValue $ws "D41" "0.006411"
Set-TextValue $ws "B42" "BKEXToken"
Set-TextValue $ws "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1053"
Set-TextValue $ws "E42" "41BKEXTokenBKK"
Set-TextValue $ws "B43" "CEJI"
Set-TextValue $ws "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.002410"
Set-TextValue $ws "E43" "42CEJICEJIWorstin24h"
Set-TextValue $ws "D44" "0.006854"
Set-TextValue $ws "E44" "43LocalTradersLCT"
Set-TextValue $ws "D45" "0.00005472"
Set-TextValue $ws "D47" "0.5503"
Set-TextValue $ws "D48" "0.006535"
Set-TextValue $ws "D49" "0.00002101"
Set-TextValue $ws "D50" "0.0002001"

Write-Output "Applied all changes"